# Pay period 26 fix: hours precision, equity/paid split, and name corrections.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 6-13: Actual Expense section ---
# The legacy "Comments" column H (which duplicated the equity hours) is retired;
# column G now holds the correct equity/paid split value, and column H is cleared.

# Row 6 - Dharam Pal
$ws.Range("F6").Value = 40
$ws.Range("G6").Value = 40
$ws.Range("H6").ClearContents()

# Row 7 - Ariful Islam
$ws.Range("F7").Value = 26.5
$ws.Range("G7").Value = 40
$ws.Range("H7").ClearContents()

# Row 8 - Raheel Shahzad
$ws.Range("F8").Value = 40
$ws.Range("G8").Value = 40
$ws.Range("H8").ClearContents()

# Row 9 - Mori Wesonga
$ws.Range("F9").Value = 40
$ws.Range("G9").Value = 40
$ws.Range("H9").ClearContents()

# Row 10 - Yulia McCoy
$ws.Range("F10").Value = 10
$ws.Range("G10").Value = 40
$ws.Range("H10").ClearContents()

# Row 11 - Pauline Nguyen
$ws.Range("F11").Value = 40
$ws.Range("G11").Value = 40
$ws.Range("H11").ClearContents()

# Row 12 - Edward Obi
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = 40
$ws.Range("H12").ClearContents()
$ws.Range("I12").Value = "Extra 11.00 hours carry over"

# Row 13 - EXPECTED EXPENSES
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").ClearContents()
$ws.Range("I13").Value = "Short 20.00 hours"

# --- Rows 22-25: Weekly Expected Expense section ---

# Row 22 - Edward Obi
$ws.Range("B22").Value = 100
$ws.Range("D22").Value = 22.5

# Row 23 - HOURS
$ws.Range("D23").Value = 40

# Row 24 - Coloring
$ws.Range("D24").Value = 40
$ws.Range("G24").Value = 10

# Row 25 - Name header row, add new multiplier cells
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 1

# --- Rows 34-36: Dev/equity reconciliation section ---

# Row 34 - Edward Obi
$ws.Range("G34").Value = 22.5

# Row 35 - name correction + equity hours
$ws.Range("A35").Value = "Dennis Fisher"
$ws.Range("G35").Value = 40

# Row 36 - name correction + equity hours
$ws.Range("A36").Value = "Forrest Cordova"
$ws.Range("G36").Value = 40
